$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the row above (row 12) into the new row 13
$ws.Range("A12:C12").Copy()
$ws.Range("A13:C13").PasteSpecial(-4122)

# Set the new row's values
$ws.Range("A13").Value2 = 45212
$ws.Range("B13").Value2 = $ws.Range("B12").Value2
$ws.Range("C13").Value2 = "Contributed technical work by aiding in resolving inconsistencies flagged by the system for employee calls"

$ws.Range("C14").Select()
